# Update the dSF (column F) values for specific rows in Sheet1.
# These rows had their "final" delta-S (dSF) value re-derived after a
# repull/push of data and a mean calculation, causing it to diverge
# from the original dS0 (column E) value it used to mirror.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    4  = -2
    5  = -2
    9  = -5
    15 = -1
    16 = -6
    18 = -1
    23 = 5
    30 = 2
    33 = -2
    36 = 0
    40 = 3
    49 = -3
    50 = 1
    52 = -3
    54 = -2
    56 = 2
    60 = -8
    63 = 1
    68 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
